# Automatic update of files.
# Updates the "Förändrad" (Changed) date column C for every data row
# from 2023-10-05 (serial 45204) to 2023-10-06 (serial 45205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
